$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 364, shifting existing rows 364:483 down to 365:484
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new data record
$ws.Cells.Item(364, 1).Value = 3
$ws.Cells.Item(364, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(364, 3).Value = "Coquimbo"
$ws.Cells.Item(364, 4).Value = 44985
$ws.Cells.Item(364, 5).Value = 5
$ws.Cells.Item(364, 6).Value = 100112012
$ws.Cells.Item(364, 7).Value = "Espinaca"
$ws.Cells.Item(364, 8).Value = "Sin especificar"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 150
$ws.Cells.Item(364, 11).Value = 6000
$ws.Cells.Item(364, 12).Value = 6500
$ws.Cells.Item(364, 13).Value = 6200
$ws.Cells.Item(364, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(364, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(364, 16).Value = 2067
$ws.Cells.Item(364, 17).Value = 3
$ws.Cells.Item(364, 18).Value = "Hortaliza"
